$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New corrected values for td_sim_1 (column C), rows 2-105
$cValues = @(
    164,
    56,
    70,
    210,
    44,
    67,
    129,
    108,
    130,
    234,
    166,
    353,
    125,
    159,
    88,
    30,
    83,
    120,
    141,
    132,
    469,
    193,
    177,
    24,
    13,
    14,
    31,
    169,
    54,
    25,
    287,
    82,
    51,
    58,
    69,
    20,
    42,
    52,
    107,
    18,
    106,
    170,
    64,
    36,
    80,
    9,
    113,
    17,
    50,
    53,
    115,
    123,
    167,
    397,
    154,
    49,
    7,
    74,
    62,
    77,
    46,
    61,
    109,
    3,
    121,
    185,
    175,
    47,
    37,
    4,
    34,
    1720,
    95,
    111,
    117,
    165,
    38,
    98,
    12,
    216,
    181,
    112,
    55,
    124,
    133,
    172,
    145,
    295,
    16,
    114,
    19,
    151,
    60,
    176,
    119,
    183,
    101,
    5,
    39,
    35,
    208,
    179,
    236,
    11
)

# New corrected values for record_atd (column D), rows 2-105 = AVERAGE(B,C)
$dValues = @(
    157,
    47,
    82,
    206.5,
    49,
    58.5,
    122.5,
    105.5,
    135.5,
    211.5,
    168,
    365.5,
    129.5,
    178,
    90,
    58,
    78,
    115.5,
    111,
    137.5,
    455.5,
    206,
    176.5,
    36,
    11,
    10,
    74.5,
    181,
    45.5,
    20.5,
    280,
    75.5,
    42.5,
    51.5,
    64.5,
    57.5,
    37,
    46.5,
    102.5,
    21,
    116.5,
    173.5,
    58.5,
    31,
    76,
    12,
    110,
    20,
    46,
    44,
    109.5,
    130.5,
    173.5,
    401,
    155,
    42.5,
    5.5,
    65.5,
    52.5,
    77.5,
    39.5,
    54,
    114.5,
    2,
    124.5,
    178,
    171.5,
    56.5,
    32,
    4,
    29.5,
    1684.5,
    94,
    106.5,
    123,
    165,
    33,
    111,
    10,
    207.5,
    182.5,
    112,
    47,
    128.5,
    135,
    172.5,
    142.5,
    297.5,
    14,
    113.5,
    58,
    153,
    54.5,
    178.5,
    120.5,
    184,
    94.5,
    3.5,
    34,
    51,
    208,
    174,
    234.5,
    10.5
)

for ($i = 0; $i -lt $cValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
}

# Summary row 106: average_simulation_TD for td_sim_1
$ws.Cells.Item(106, 3).Value = 122.5961538461538
